# Re-shuffle the observation rows 19-28 and 31-35: each row keeps its own
# row number, but the full record (columns A:AY) that used to live in one
# row now lives in another row, per the mapping below (derived from the
# target OOXML diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (i.e. target row receives source row's old content)
$map = @{
    19 = 20
    20 = 24
    21 = 23
    22 = 21
    23 = 25
    24 = 19
    25 = 26
    26 = 22
    27 = 28
    28 = 27
    31 = 35
    32 = 33
    33 = 31
    34 = 32
    35 = 34
}

$firstCol = 1   # A
$lastCol  = 51  # AY

# Columns that hold genuine numbers (everything else on this sheet is
# stored as text, even when the text looks like a number, e.g. the
# "Antal" column I holding "4").
$numericCols = @(1, 2, 5, 17, 18, 19)       # A, B, E, Q, R, S
# Columns that hold booleans.
$boolCols = @(30, 31, 33)                   # AD, AE, AG
# Text columns whose values can look like numbers/dates and therefore need
# an explicit text format so Value2 doesn't get auto-parsed into a date
# serial / number on write.
$forceTextCols = @(9, 25, 26, 27, 28)       # I, Y, Z, AA, AB

function Values-Equal($a, $b) {
    if ($a -is [bool] -or $b -is [bool]) {
        if ($a -is [bool] -and $b -is [bool]) { return $a -eq $b }
        return $false
    }
    if ($null -eq $a) { $a = "" }
    if ($null -eq $b) { $b = "" }
    return ([string]$a) -ceq ([string]$b)
}

# Snapshot every row that participates (as source or destination) BEFORE
# any writes happen, so overlapping/cyclic moves don't clobber data we
# still need to read, and so we can diff old-vs-new per cell.
$rowsInvolved = @{}
foreach ($k in $map.Keys) { $rowsInvolved[$k] = $true }
foreach ($v in $map.Values) { $rowsInvolved[$v] = $true }

$snapshots = @{}
foreach ($r in $rowsInvolved.Keys) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshots[$r] = $rowVals
}

foreach ($dstRow in ($map.Keys | Sort-Object)) {
    $srcRow = $map[$dstRow]
    $oldVals = $snapshots[$dstRow]
    $newVals = $snapshots[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $oldVal = $oldVals[$c]
        $newVal = $newVals[$c]
        if (Values-Equal $oldVal $newVal) {
            continue
        }
        $cell = $ws.Cells.Item($dstRow, $c)
        $isBlank = ($null -eq $newVal -or $newVal -eq "")
        if ($boolCols -contains $c) {
            if ($isBlank) {
                $cell.ClearContents()
            } else {
                $cell.Value2 = [bool]$newVal
            }
        } elseif ($numericCols -contains $c) {
            if ($isBlank) {
                $cell.ClearContents()
            } else {
                $cell.Value2 = $newVal
            }
        } else {
            if ($forceTextCols -contains $c) {
                $cell.NumberFormat = "@"
            }
            if ($isBlank) {
                $cell.ClearContents()
            } else {
                $cell.Value2 = [string]$newVal
            }
        }
    }
}

Write-Output "done"
